$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Starting point (paragraph 3, the last paragraph in the doc):
#   "He tries to escape them by running to the washroom"   + _GoBack bookmark
#    (bookmark is collapsed, sitting right after the text, before the
#     paragraph mark)
#
# Target:
#   Paragraph 3: "He tries to escape them by running to the washroom" run,
#                plus a *separate* run containing just "."
#   Paragraph 4 (new): "He runs into the cubicle, closes the door and sits
#                down…that’s when it happened." with the _GoBack
#                bookmark now collapsed at its end.
# ---------------------------------------------------------------------

$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range
$r.MoveEnd(1, -1) | Out-Null
$endOfWashroom = $r.End

# Drop the old _GoBack bookmark -- Word re-anchors "last edit position" to
# wherever the new text lands, so we will recreate it at the new location.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Add the closing "." as a run of its own -----------------------------
# A plain InsertAfter would just extend/merge into the existing run since
# both share identical (default) formatting, so instead we park the period
# in a throw-away paragraph and then splice the intervening paragraph mark
# back out; that leaves the two runs distinct.
$r.InsertParagraphAfter()
$periodPara = $d.Paragraphs.Last
$rPeriod = $periodPara.Range
$rPeriod.MoveEnd(1, -1) | Out-Null
$rPeriod.InsertAfter(".")
$d.Range($endOfWashroom, $endOfWashroom + 1).Delete()

# --- Start the new story paragraph ---------------------------------------
$lastPara2 = $d.Paragraphs.Last
$r2 = $lastPara2.Range
$r2.MoveEnd(1, -1) | Out-Null
$r2.InsertParagraphAfter()

$storyPara = $d.Paragraphs.Last
$r3 = $storyPara.Range
$r3.MoveEnd(1, -1) | Out-Null
$r3.InsertAfter("He runs into the cubicle, closes the door and sits down…that’s when it happened.")

# --- Re-anchor _GoBack, collapsed, at the very end of the new paragraph --
$finalPara = $d.Paragraphs.Last
$rEnd = $finalPara.Range
$rEnd.MoveEnd(1, -1) | Out-Null
$storyEnd = $rEnd.End

# A Bookmarks.Add() with a zero-length range exactly at a paragraph's text
# end gets mis-resolved, so nudge it: insert a throw-away character after
# the intended bookmark position, bookmark the still-collapsed point in
# front of it, then delete the throw-away character again.
$d.Range($storyEnd, $storyEnd).InsertAfter("X")
$bmRange = $d.Range($storyEnd, $storyEnd)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
$d.Range($storyEnd, $storyEnd + 1).Delete()
